$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.567.89"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +12.80%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'1.841.51"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +9.85%  "
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'232.33"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +5.72%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'0.571"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +7.73%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'31.84"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  +7.61%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'47.19"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +6.34%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.290"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +9.85%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.0687"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  +6.92%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.0935"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +3.29%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'2.110.30"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +10.06%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'1.848.26"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +10.37%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.661"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +8.44%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'34.548.39"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +12.70%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'10.37"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +3.71%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'4.34"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +7.88%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'70.81"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  +6.75%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'260.16"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  +7.20%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'0.0₃0766"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +6.22%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'0.998"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  -0.16%  "
$ws.Range("E22").ClearFormats()
$ws.Range("E23").Value = "'  +7.36%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'4.42"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +3.91%  "
$ws.Range("E24").ClearFormats()
$ws.Range("E25").Value = "'  +4.31%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'159.22"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  +0.18%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'16.93"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  +6.92%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'0.118"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  +5.00%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'7.24"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  +8.26%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  -0.04%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'3.96"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  +14.36%  "
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = "'  +7.49%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'1.22"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +7.00%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'3.63"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +10.19%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'1.564.26"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  +4.19%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'1.83"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +2.31%  "
$ws.Range("E36").ClearFormats()
$ws.Range("E37").Value = "'  +6.62%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'0.650"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  +8.57%  "
$ws.Range("E38").ClearFormats()
$ws.Range("B39").Value = "'VeChain"
$ws.Range("B39").ClearFormats()
$ws.Range("C39").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C39").ClearFormats()
$ws.Range("D39").Value = "'0.0193"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +7.94%  "
$ws.Range("E39").ClearFormats()
$ws.Range("B40").Value = "'Aave"
$ws.Range("B40").ClearFormats()
$ws.Range("C40").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C40").ClearFormats()
$ws.Range("D40").Value = "'85.99"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  +2.48%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'2.84"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  +6.11%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.926"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +10.28%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'2.35"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +2.55%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'2.16"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +6.43%  "
$ws.Range("E44").ClearFormats()
$ws.Range("B45").Value = "'Kaspa"
$ws.Range("B45").ClearFormats()
$ws.Range("C45").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C45").ClearFormats()
$ws.Range("D45").Value = "'0.0529"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +5.64%  "
$ws.Range("E45").ClearFormats()
$ws.Range("B46").Value = "'MinaProtocolToken"
$ws.Range("B46").ClearFormats()
$ws.Range("C46").Value = "'https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina"
$ws.Range("C46").ClearFormats()
$ws.Range("D46").Value = "'1.06"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +152.86%  "
$ws.Range("E46").ClearFormats()
$ws.Range("E47").Value = "'  +6.17%  "
$ws.Range("E47").ClearFormats()
$ws.Range("B48").Value = "'RocketPoolETH"
$ws.Range("B48").ClearFormats()
$ws.Range("C48").Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("C48").ClearFormats()
$ws.Range("D48").Value = "'1.998.73"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +10.27%  "
$ws.Range("E48").ClearFormats()
$ws.Range("B49").Value = "'InjectiveProtocol"
$ws.Range("B49").ClearFormats()
$ws.Range("C49").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C49").ClearFormats()
$ws.Range("D49").Value = "'12.44"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +19.28%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'5.87"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  +5.74%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'1.00"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  -0.01%  "
$ws.Range("E51").ClearFormats()
